# Add support for timevals: convert the "Data" sheet's year column (C)
# from plain numeric years (2018..2021) to quarter-style text labels
# (2000Q1..2000Q4), and move the active-sheet/selection focus from the
# "Codelists" sheet to the "Data" sheet.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Data")
$codelistsSheet = $wb.Worksheets.Item("Codelists")

# Map the old numeric year values to the new quarter text labels.
$yearToQuarter = @{
    2018 = "2000Q1"
    2019 = "2000Q2"
    2020 = "2000Q3"
    2021 = "2000Q4"
}

# Column C holds the year values for rows 2 through 85.
for ($row = 2; $row -le 85; $row++) {
    $cell = $dataSheet.Cells.Item($row, 3)
    $year = [int]$cell.Value()
    $cell.Value = $yearToQuarter[$year]
}

# The "Codelists" sheet was previously the selected/active tab; move that
# to the "Data" sheet instead, along with an updated selection.
[void]$codelistsSheet.Select()
[void]$codelistsSheet.Range("C11").Select()

[void]$dataSheet.Select()
[void]$dataSheet.Range("E16").Select()
